$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2..44 (Application, Year, Count)
$data = @(
    @("Bio-implementation", 2004, 1),
    @("Bio-implementation", 2013, 2),
    @("Bio-implementation", 2014, 1),
    @("Bio-implementation", 2016, 1),
    @("Bio-implementation", 2017, 1),
    @("Bio-implementation", 2021, 2),
    @("Bio-implementation", 2022, 4),
    @("Bio-implementation", 2023, 2),
    @("Bio-implementation", 2024, 3),
    @("Bio-implementation", 2025, 2),
    @("Dataset generation", 2004, 1),
    @("Dataset generation", 2005, 1),
    @("Dataset generation", 2010, 1),
    @("Dataset generation", 2011, 1),
    @("Dataset generation", 2012, 2),
    @("Dataset generation", 2013, 2),
    @("Dataset generation", 2014, 1),
    @("Dataset generation", 2016, 1),
    @("Dataset generation", 2017, 4),
    @("Dataset generation", 2018, 2),
    @("Dataset generation", 2019, 4),
    @("Dataset generation", 2020, 5),
    @("Dataset generation", 2021, 1),
    @("Dataset generation", 2022, 6),
    @("Dataset generation", 2023, 8),
    @("Dataset generation", 2024, 4),
    @("Explainability", 2021, 1),
    @("Explainability", 2023, 2),
    @("Explainability", 2024, 1),
    @("Explainability", 2025, 1),
    @("MAC layer", 2024, 1),
    @("PHY layer", 2017, 2),
    @("PHY layer", 2018, 7),
    @("PHY layer", 2019, 2),
    @("PHY layer", 2020, 4),
    @("PHY layer", 2021, 5),
    @("PHY layer", 2022, 5),
    @("PHY layer", 2023, 11),
    @("PHY layer", 2024, 10),
    @("PHY layer", 2025, 3),
    @("Upper layers", 2017, 1),
    @("Upper layers", 2020, 3),
    @("Upper layers", 2024, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
}

$wb.Save()
